# "Generate Report for Handoff"
# Updates the localization-status report: the "In Translation" status
# becomes "Ready for handoff" and the associated generation timestamps
# advance by one minute, now that the handoff package has been produced.
# Column widths on the affected "Status" columns are widened so the new,
# longer status text is no longer truncated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-14 16:52:39"
$overview.Columns.Item(5).ColumnWidth = 16.333333
$overview.Columns.Item(6).ColumnWidth = 16.333333

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-14 16:52:31"
$zhcn.Columns.Item(3).ColumnWidth = 16.333333

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-14 16:52:39"
$dede.Columns.Item(3).ColumnWidth = 16.333333
